# "selected module to work on"
# Fill in the Author column for several module rows in the HRM Modules
# table with the name of the developer who picked up that module.
#
# Table layout: column 1 = Module name, column 2 = Author.
# A handful of rows whose Author cell only contained a bare "/" (meaning
# "no one yet") get a name typed in front of the slash; one row whose
# Author cell was completely empty gets just the name typed in.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Insert-AuthorName($table, $row, $name) {
    $cell = $table.Cell($row, 2)
    $r = $cell.Range
    $r.Collapse(1)
    $r.InsertBefore($name)
}

# Patient module
Insert-AuthorName $t 3  "Jide"   # Add New Patient
Insert-AuthorName $t 4  "Jide"   # View Patient List {Edit, View & Delete}

# Nurse module
Insert-AuthorName $t 13 "jide"   # Add New Nurse
Insert-AuthorName $t 14 "jide"   # View Nurse List {Edit, View & Delete}
Insert-AuthorName $t 15 "Jide"   # Attendance Page

# Doctor module
Insert-AuthorName $t 19 "Jide"   # Add New Doctor (cell was empty)
Insert-AuthorName $t 20 "Jide"   # View Doctor List {Edit, View & Delete}
Insert-AuthorName $t 21 "jide"   # Attendance Page
